$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on numeric-looking price cells so Excel does not
# reinterpret them as numbers (they are text in the source data).
$textCells = @("D5", "D6", "D8", "D9", "D10", "D11", "D13", "D14", "D15", "D16", "D18", "D22", "D23", "D25", "D27", "D28", "D29", "D30", "D31", "D32", "D34", "D36", "D37", "D39", "D40", "D41", "D42", "D43", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($c in $textCells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range("D2").Value = '29.485.58'
$ws.Range("E2").Value = '  +0.12%  '
$ws.Range("D3").Value = '1.857.67'
$ws.Range("E3").Value = '  +0.44%  '
$ws.Range("E4").Value = '  +0.17%  '
$ws.Range("D5").Value = '241.75'
$ws.Range("E5").Value = '  +0.32%  '
$ws.Range("D6").Value = '0.6334'
$ws.Range("E6").Value = '  +0.99%  '
$ws.Range("E7").Value = '  +0.13%  '
$ws.Range("D8").Value = '0.07584'
$ws.Range("E8").Value = '  -1.16%  '
$ws.Range("D9").Value = '0.2929'
$ws.Range("E9").Value = '  +0.48%  '
$ws.Range("D10").Value = '24.63'
$ws.Range("E10").Value = '  -0.82%  '
$ws.Range("D11").Value = '0.07766'
$ws.Range("E11").Value = '  +0.17%  '
$ws.Range("D12").Value = '1.857.82'
$ws.Range("E12").Value = '  +0.33%  '
$ws.Range("D13").Value = '5.046'
$ws.Range("E13").Value = '  +0.26%  '
$ws.Range("D14").Value = '0.6860'
$ws.Range("E14").Value = '  +0.69%  '
$ws.Range("D15").Value = '0.00001047'
$ws.Range("E15").Value = '  -2.32%  '
$ws.Range("D16").Value = '83.50'
$ws.Range("E16").Value = '  -0.01%  '
$ws.Range("D17").Value = '2.116.69'
$ws.Range("E17").Value = '  -0.71%  '
$ws.Range("D18").Value = '6.159'
$ws.Range("E18").Value = '  -0.21%  '
$ws.Range("D19").Value = '29.493.70'
$ws.Range("E19").Value = '  +0.10%  '
$ws.Range("E20").Value = '  +0.78%  '
$ws.Range("D22").Value = '1.003'
$ws.Range("E22").Value = '  +0.17%  '
$ws.Range("D23").Value = '7.523'
$ws.Range("E23").Value = '  +1.27%  '
$ws.Range("E24").Value = '  +0.17%  '
$ws.Range("D25").Value = '159.61'
$ws.Range("E25").Value = '  +0.96%  '
$ws.Range("E26").Value = '  +2.03%  '
$ws.Range("D27").Value = '8.474'
$ws.Range("E27").Value = '  +0.82%  '
$ws.Range("D28").Value = '17.75'
$ws.Range("E28").Value = '  +0.29%  '
$ws.Range("D29").Value = '1.423'
$ws.Range("E29").Value = '  +5.14%  '
$ws.Range("D30").Value = '1.482'
$ws.Range("E30").Value = '  +0.94%  '
$ws.Range("D31").Value = '0.05709'
$ws.Range("E31").Value = '  +1.10%  '
$ws.Range("D32").Value = '4.159'
$ws.Range("E32").Value = '  +0.96%  '
$ws.Range("E33").Value = '  +0.80%  '
$ws.Range("D34").Value = '1.835'
$ws.Range("E34").Value = '  -0.43%  '
$ws.Range("E35").Value = '  -0.15%  '
$ws.Range("D36").Value = '0.6955'
$ws.Range("E36").Value = '  -1.78%  '
$ws.Range("D37").Value = '2.597'
$ws.Range("E37").Value = '  +0.10%  '
$ws.Range("D38").Value = '1.257.10'
$ws.Range("E38").Value = '  +2.46%  '
$ws.Range("D39").Value = '0.01835'
$ws.Range("E39").Value = '  +2.51%  '
$ws.Range("D40").Value = '2.785'
$ws.Range("E40").Value = '  +0.68%  '
$ws.Range("D41").Value = '6.522'
$ws.Range("E41").Value = '  -0.37%  '
$ws.Range("D42").Value = '0.9085'
$ws.Range("E42").Value = '  +0.59%  '
$ws.Range("D43").Value = '1.003'
$ws.Range("E43").Value = '  +0.11%  '
$ws.Range("D44").Value = '2.019.62'
$ws.Range("E44").Value = '  -0.92%  '
$ws.Range("D45").Value = '101.77'
$ws.Range("E45").Value = '  -0.05%  '
$ws.Range("D46").Value = '66.28'
$ws.Range("E46").Value = '  +0.41%  '
$ws.Range("D47").Value = '7.170'
$ws.Range("E47").Value = '  +0.07%  '
$ws.Range("D48").Value = '0.1172'
$ws.Range("E48").Value = '  +1.45%  '
$ws.Range("D49").Value = '9.041'
$ws.Range("E49").Value = '  +0.17%  '
$ws.Range("D50").Value = '0.3983'
$ws.Range("E50").Value = '  -0.79%  '
$ws.Range("D51").Value = '1.680'
$ws.Range("E51").Value = '  +0.45%  '
